$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for changed rows.
# NumberFormat "@" forces the Price cells to stay text (matches the
# source data, which stores prices as literal strings, some of which
# look numeric e.g. "197.90").
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "76.362.88"
$ws.Cells.Item(2, 5).Value = "  -0.06%  "
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "3.036.73"
$ws.Cells.Item(3, 5).Value = "  +3.38%  "
$ws.Cells.Item(4, 5).Value = "  +0.11%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "197.90"
$ws.Cells.Item(5, 5).Value = "  -1.78%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "616.83"
$ws.Cells.Item(6, 5).Value = "  +2.95%  "
$ws.Cells.Item(7, 5).Value = "  +0.10%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.547"
$ws.Cells.Item(8, 5).Value = "  -1.30%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.205"
$ws.Cells.Item(9, 5).Value = "  +4.06%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "3.035.19"
$ws.Cells.Item(10, 5).Value = "  +3.68%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.434"
$ws.Cells.Item(11, 5).Value = "  -1.41%  "
$ws.Cells.Item(12, 5).Value = "  -0.73%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "5.22"
$ws.Cells.Item(13, 5).Value = "  +5.15%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "3.595.19"
$ws.Cells.Item(14, 5).Value = "  +3.67%  "
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "28.71"
$ws.Cells.Item(15, 5).Value = "  +1.52%  "
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "75.876.51"
$ws.Cells.Item(16, 5).Value = "  -0.42%  "
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "0.0000192"
$ws.Cells.Item(17, 5).Value = "  +1.12%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "3.034.25"
$ws.Cells.Item(18, 5).Value = "  +3.83%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "13.52"
$ws.Cells.Item(19, 5).Value = "  +1.08%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "8.97"
$ws.Cells.Item(20, 5).Value = "  +2.41%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "379.22"
$ws.Cells.Item(21, 5).Value = "  +1.33%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "2.34"
$ws.Cells.Item(22, 5).Value = "  +1.09%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "4.35"
$ws.Cells.Item(23, 5).Value = "  -0.14%  "
$ws.Cells.Item(24, 5).Value = "  +3.55%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "72.53"
$ws.Cells.Item(25, 5).Value = "  +0.26%  "
$ws.Cells.Item(26, 5).Value = "  +0.03%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "4.32"
$ws.Cells.Item(27, 5).Value = "  -0.30%  "
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "9.69"
$ws.Cells.Item(28, 5).Value = "  -0.53%  "
$ws.Cells.Item(29, 5).Value = "  -1.21%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "0.999"
$ws.Cells.Item(30, 5).Value = "  -0.19%  "
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "8.23"
$ws.Cells.Item(31, 5).Value = "  +3.57%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "1.39"
$ws.Cells.Item(32, 5).Value = "  +0.32%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "492.73"
$ws.Cells.Item(33, 5).Value = "  -1.89%  "
$ws.Cells.Item(34, 5).Value = "  +3.89%  "
$ws.Cells.Item(35, 5).Value = "  +0.10%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "20.53"
$ws.Cells.Item(36, 5).Value = "  +1.02%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "163.19"
$ws.Cells.Item(37, 5).Value = "  -1.57%  "
$ws.Cells.Item(38, 5).Value = "  +1.95%  "
$ws.Cells.Item(39, 5).Value = "  +4.25%  "
$ws.Cells.Item(42, 5).Value = "  -3.44%  "
$ws.Cells.Item(43, 5).Value = "  +0.04%  "
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.792"
$ws.Cells.Item(44, 5).Value = "  +20.00%  "
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "5.05"
$ws.Cells.Item(45, 5).Value = "  +1.09%  "
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "41.78"
$ws.Cells.Item(46, 5).Value = "  +4.07%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "1.25"
$ws.Cells.Item(47, 5).Value = "  +3.89%  "
$ws.Cells.Item(48, 5).Value = "  -1.85%  "
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "2.40"
$ws.Cells.Item(49, 5).Value = "  +2.41%  "
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.598"
$ws.Cells.Item(50, 5).Value = "  +2.32%  "
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "3.87"
$ws.Cells.Item(51, 5).Value = "  -0.32%  "

# Rows 40-41: Aave overtakes PolygonEcosystemToken, so they swap
# positions in the ranking; update Coin, Link, Price and Volume.
$ws.Cells.Item(40, 2).Value = "Aave"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "191.26"
$ws.Cells.Item(40, 5).Value = "  +6.94%  "
$ws.Cells.Item(41, 2).Value = "PolygonEcosystemToken"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.377"
$ws.Cells.Item(41, 5).Value = "  -2.07%  "
